$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Стикеры")

# Clear the old contents of the used range first
$ws.Range("A1:D4").Clear()

# New header row
$ws.Range("B1").Value = "emoji"
$ws.Range("C1").Value = "file_id"
$ws.Range("D1").Value = "sticker_id"
$ws.Range("E1").Value = "key_word"

# Data rows
$ws.Range("D2").Value = "CAACAgIAAxkBAANBYFnKfcWD9t6m_8-4LD8clr4e4wcAAlkAAwr8wgXOU7sZfH5zGx4E"
$ws.Range("E2").Value = "смешная шутка"

$ws.Range("D3").Value = "CAACAgIAAxkBAANEYFnLQxjtNi5MTuMghLPi9mJjD3MAAg0EAALPX4sHmuYS8a7yxGQeBA"
$ws.Range("E3").Value = "я не хочу брать Иерусалим"

# Apply font formatting to D2:D3 (Tahoma 10, black)
$fmtRange = $ws.Range("D2:D3")
$fmtRange.Font.Name = "Tahoma"
$fmtRange.Font.Size = 10
$fmtRange.Font.Color = 0

# Column widths
$ws.Columns.Item(4).ColumnWidth = 29.42578125
$ws.Columns.Item(5).ColumnWidth = 10.28515625

# Selection
$ws.Range("E6").Select()
